# GUI: Renamed the suites.
# The "SourceAttributes" test suite is split into two suites:
#   - SourceAttributeFields (0 automated, 3 total) - status "Suited to Manual"
#   - SourceAttributeLists  (0 automated, 11 total) - status "Suited to Manual"
# This inserts a new row for the table on Sheet1 and updates the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 7; this shifts old rows 7-9 down to 8-10 and keeps all
# formulas ($A$3:A38 etc.) correctly re-pointed to the new extents.
$ws.Rows("7:7").Insert()

# New row 7: SourceAttributeFields
$ws.Range("A7").Value = "SourceAttributeFields"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Suited to Manual"

# Row 8 used to hold the old "SourceAttributes" data (now shifted from row 7);
# rename/update it to the new "SourceAttributeLists" suite.
$ws.Range("A8").Value = "SourceAttributeLists"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = "Suited to Manual"

# Rows 9 and 10 (old rows 8 and 9 - Step3Fields / WorkspaceFields) keep their
# original values; Insert already shifted them down so nothing else to do.

# Extend the conditional formatting range that covered D3:D51 by one row so
# it still reaches the bottom of the (now one-row-taller) table.
$fcs = $ws.Range("D3").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D3:D52"))

# Move the active selection to D2, matching the saved selection state.
$ws.Range("D2").Select() | Out-Null
